# Update the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions scheduled refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold values that look numeric
# (e.g. "0.513", "6.54", "  -0.08%  ") even though the sheet stores them
# as plain text. Force the whole data range to Text format before
# writing so Excel does not silently convert them to real numbers, then
# restore the default "Normal" style afterwards so no stray
# number-formatting is left behind on the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '67.698.77'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '3.772.77'
$ws.Range("E3").Value = '  -0.73%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '598.57'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = '163.09'
$ws.Range("E6").Value = '  -2.49%  '
$ws.Range("D7").Value = '3.769.17'
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("E10").Value = '  -2.91%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").Value = '6.54'
$ws.Range("E12").Value = '  +3.92%  '
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("D15").Value = '4.403.31'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '3.789.10'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '67.719.14'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").Value = '18.23'
$ws.Range("E18").Value = '  -1.78%  '
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").Value = '6.99'
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = '457.01'
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("E22").Value = '  -4.39%  '
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").Value = '82.71'
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  -6.27%  '
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '9.84'
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").Value = '3.916.02'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("E33").Value = '  -6.77%  '
$ws.Range("D34").Value = '28.82'
$ws.Range("E34").Value = '  -2.46%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").Value = '0.0988'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("E38").Value = '  +2.16%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").Value = '0.978'
$ws.Range("E40").Value = '  -2.19%  '

# Rows 41/42 swapped ranking order: dogwifhat now sits above
# FirstDigitalUSD, so the coin name/link/price/volume cells are
# rewritten in place.
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '3.15'
$ws.Range("E41").Value = '  -5.97%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D44").Value = '47.18'
$ws.Range("E44").Value = '  -1.92%  '
$ws.Range("D45").Value = '42.99'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").Value = '152.21'
$ws.Range("E46").Value = '  +3.27%  '
$ws.Range("E47").Value = '  -2.37%  '
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("D51").Value = '384.88'
$ws.Range("E51").Value = '  -2.55%  '

# Restore the original (default) cell style now that the text values
# are safely in place.
$dataRange.Style = "Normal"
